$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Electric" living-expense line is being inserted above the existing
# "Other" / {other_living_expenses} line (row 13), pushing that pair down to
# row 14. No sheet row is physically inserted -- just the cell contents of
# rows 13/14 shift -- so move the old row-13 content down to row 14 first,
# then overwrite row 13 with the new "Electric" line.

$ws.Range("A14").Value = $ws.Range("A13").Value()
$ws.Range("B14").Value = $ws.Range("B13").Value()

$ws.Range("A13").Value = "Electric"
# Empty placeholder value for the new line -- stored as a (quote-prefixed)
# text cell rather than a truly blank cell.
$ws.Range("B13").Value = "'"

# Row 15 picks up the uniform 19.5pt row height used by the rest of the sheet.
$ws.Rows(15).RowHeight = 19.5
